$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# like "104.70" or "0.160" keep their exact text representation
# instead of being parsed into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "72.116.77"
$ws.Range("E2").Value = "  +4.11%  "
$ws.Range("D3").Value = "4.035.53"
$ws.Range("E3").Value = "  +3.63%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "520.91"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").Value = "147.33"
$ws.Range("E6").Value = "  +2.02%  "
$ws.Range("D7").Value = "0.724"
$ws.Range("E7").Value = "  +18.64%  "
$ws.Range("D8").Value = "4.026.88"
$ws.Range("E8").Value = "  +3.58%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").Value = "0.774"
$ws.Range("E10").Value = "  +7.87%  "
$ws.Range("D11").Value = "0.175"
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").Value = "48.54"
$ws.Range("E13").Value = "  +15.51%  "
$ws.Range("D14").Value = "11.12"
$ws.Range("E14").Value = "  +8.57%  "
$ws.Range("D15").Value = "4.678.38"
$ws.Range("D16").Value = "4.051.39"
$ws.Range("E16").Value = "  +3.60%  "
$ws.Range("D17").Value = "21.26"
$ws.Range("E17").Value = "  +7.68%  "
$ws.Range("D18").Value = "14.24"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("D19").Value = "1.21"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "72.136.18"
$ws.Range("E21").Value = "  +4.20%  "
$ws.Range("D22").Value = "443.93"
$ws.Range("E22").Value = "  +4.68%  "
$ws.Range("D23").Value = "104.70"
$ws.Range("E23").Value = "  +19.21%  "
$ws.Range("D24").Value = "3.57"
$ws.Range("E24").Value = "  +5.48%  "
$ws.Range("D25").Value = "15.05"
$ws.Range("E25").Value = "  +6.59%  "
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").Value = "11.55"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("D28").Value = "11.05"
$ws.Range("E28").Value = "  +4.83%  "
$ws.Range("D29").Value = "37.86"
$ws.Range("E29").Value = "  +4.17%  "
$ws.Range("E30").Value = "  +2.45%  "
$ws.Range("D31").Value = "3.25"
$ws.Range("E31").Value = "  +14.58%  "
$ws.Range("D32").Value = "13.73"
$ws.Range("E32").Value = "  +4.24%  "
$ws.Range("E33").Value = "  +3.47%  "
$ws.Range("D34").Value = "675.94"
$ws.Range("E34").Value = "  -1.70%  "
$ws.Range("D35").Value = "6.76"
$ws.Range("E35").Value = "  +14.64%  "
$ws.Range("D36").Value = "67.12"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("D37").Value = "42.48"
$ws.Range("E37").Value = "  +6.60%  "
$ws.Range("D38").Value = "0.0₃0865"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").Value = "0.426"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").Value = "3.51"
$ws.Range("E40").Value = "  +5.57%  "
$ws.Range("D41").Value = "0.153"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "0.0502"
$ws.Range("E43").Value = "  +3.90%  "
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").Value = "0.160"
$ws.Range("E46").Value = "  +14.02%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "2.71"
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "3.48"
$ws.Range("E48").Value = "  +1.88%  "
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").Value = "9.28"
$ws.Range("E50").Value = "  +8.47%  "
$ws.Range("E51").Value = "  -1.16%  "
